$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header row (row 1) ---
# Existing header cells: A1="Поз." (keep), B1="Наименование", C1="Обозначение"
# Target header cells:    A1="Поз.", B1="Обозначение", C1="Наименование",
#                          D1="Объем ед. м3", E1="Примечание"

# Copy the existing header style onto the two brand-new header cells first,
# so they end up on the same style index as the rest of the header row.
$ws.Range("A1").Copy()
$ws.Range("D1").PasteSpecial(-4122)
$ws.Range("E1").PasteSpecial(-4122)
$excel.CutCopyMode = $false

$ws.Range("A1").Value = "Поз."
$ws.Range("B1").Value = "Обозначение"
$ws.Range("C1").Value = "Наименование"
$ws.Range("D1").Value = "Объем ед. м3"
$ws.Range("E1").Value = "Примечание"

# --- Data rows ---
# Remove old rows 2-4 content, then write the single new data row (row 2).
$ws.Range("A2:C4").ClearContents()

# A2, D2 and E2 hold digit-only / numeric-looking text, so force them to be
# stored as text (not auto-converted to numbers) by pre-formatting as Text,
# then restoring the plain (unstyled) format so they don't pick up a new
# style index.
$ws.Range("A2").NumberFormat = "@"
$ws.Range("A2").Value = "2"

$ws.Range("B2").Value = "СТБ 1437"
$ws.Range("C2").Value = "ППТ-15-А-Р-2000x620x180"

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "1"

$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = "0.22319999999999998"

# Reset those three cells back to the default (unstyled) format so they
# match the rest of the data row.
$ws.Range("B2").Copy()
$ws.Range("A2").PasteSpecial(-4122)
$ws.Range("D2").PasteSpecial(-4122)
$ws.Range("E2").PasteSpecial(-4122)
$excel.CutCopyMode = $false
